$wb = $excel.ActiveWorkbook

$wsTime  = $wb.Worksheets.Item("time log")
$wsTasks = $wb.Worksheets.Item("tasks")
$wsLinks = $wb.Worksheets.Item("useful links")

# ---------------------------------------------------------------------------
# "useful links" sheet: two new rows of reference material (rows 12-13)
# Added first so the shared-string table gets the same allocation order as
# the authoring session (links before the "time log" task description).
# ---------------------------------------------------------------------------
$wsLinks.Range("D12").Value = "https://www.youtube.com/watch?v=M9lZ5Qy5S2s"
$wsLinks.Range("E12").Value = "i2c interface, arduino mpu6050"

# ---------------------------------------------------------------------------
# "time log" sheet: two new time-tracking entries (rows 21-22)
# ---------------------------------------------------------------------------
$wsTime.Range("E20:F20").Copy()
$wsTime.Range("E21").PasteSpecial(-4122)  # xlPasteFormats
$wsTime.Range("I20").Copy()
$wsTime.Range("I21").PasteSpecial(-4122)

$wsTime.Range("E21").Value = 43398
$wsTime.Range("F21").Value = "read docs & read data from imu on kl25z"
$wsTime.Range("I21").Value = 1
$wsTime.Rows.Item(21).RowHeight = 29

$wsTime.Range("E20:F20").Copy()
$wsTime.Range("E22").PasteSpecial(-4122)  # xlPasteFormats
$wsTime.Range("I20").Copy()
$wsTime.Range("I22").PasteSpecial(-4122)

$wsTime.Range("E22").Value = 43402
$wsTime.Range("F22").Value = "connect motor shield to kl25z and motor battery"
$wsTime.Range("I22").Value = 2
$wsTime.Rows.Item(22).RowHeight = 29

# ---------------------------------------------------------------------------
# back to "useful links" for the second new row + column D widening
# ---------------------------------------------------------------------------
$wsLinks.Range("D13").Value = "https://www.instructables.com/id/Arduino-Nano-Control-Brushelss-DC-Motor-With-L9110/"
$wsLinks.Range("E13").Value = "connect motor controller (test) to board"

$wsLinks.Columns.Item(4).ColumnWidth = 77

# ---------------------------------------------------------------------------
# Selection / view-state updates recorded by the diff. Selecting on a sheet
# activates it, so we select on each sheet in turn and finish by re-activating
# "useful links" (the tab that is actually active in the saved workbook).
# ---------------------------------------------------------------------------
$wsTime.Range("H20").Select()
$wsTasks.Range("D8").Select()
$wsLinks.Range("E14").Select()
$wsLinks.Activate()
